$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 167 (id 165)
$ws.Range("B167").Value = 6992693
$ws.Range("E167").Value = 45382.3125
$ws.Range("F167").Value = "Buriram United"
$ws.Range("G167").Value = "Chiangrai Utd"
$ws.Range("K167").Value = 1.4
$ws.Range("L167").Value = 4.333
$ws.Range("M167").Value = 6
$ws.Range("N167").Value = 1.363
$ws.Range("O167").Value = 4.5
$ws.Range("P167").Value = 6.5
$ws.Range("Q167").Value = -1.5
$ws.Range("R167").Value = 1.95
$ws.Range("S167").Value = 1.85
$ws.Range("T167").Value = 3.25
$ws.Range("U167").Value = 1.975
$ws.Range("V167").Value = 1.825

# Row 168 (id 166)
$ws.Range("B168").Value = 6992335
$ws.Range("E168").Value = 45382.33333333334
$ws.Range("F168").Value = "Muang Thong United"
$ws.Range("G168").Value = "Police Tero FC"
$ws.Range("K168").Value = 1.444
$ws.Range("L168").Value = 4.2
$ws.Range("M168").Value = 5.5
$ws.Range("N168").Value = 1.4
$ws.Range("O168").Value = 4.333
$ws.Range("P168").Value = 5.75
$ws.Range("Q168").Value = -1.25
$ws.Range("R168").Value = 1.825
$ws.Range("S168").Value = 1.975
$ws.Range("T168").Value = 3.25
$ws.Range("U168").Value = 1.9
$ws.Range("V168").Value = 1.9

# Row 169 (id 167)
$ws.Range("B169").Value = 6992692
$ws.Range("E169").Value = 45382.35416666666
$ws.Range("F169").Value = "Bangkok United"
$ws.Range("G169").Value = "Chonburi"
$ws.Range("K169").Value = 1.444
$ws.Range("L169").Value = 4.2
$ws.Range("M169").Value = 5.5
$ws.Range("N169").Value = 1.4
$ws.Range("O169").Value = 4.2
$ws.Range("P169").Value = 6
$ws.Range("Q169").Value = -1.25

# Row 170 (id 168)
$ws.Range("B170").Value = 6992688
$ws.Range("E170").Value = 45382.375
$ws.Range("F170").Value = "Khonkaen United"
$ws.Range("G170").Value = "Nakhon Pathom FC"
$ws.Range("K170").Value = 2.7
$ws.Range("L170").Value = 3.6
$ws.Range("M170").Value = 2.15
$ws.Range("N170").Value = 2.5
$ws.Range("O170").Value = 3.4
$ws.Range("P170").Value = 2.4
$ws.Range("Q170").Value = 0
$ws.Range("R170").Value = 1.9
$ws.Range("S170").Value = 1.9
$ws.Range("T170").Value = 2.75
$ws.Range("U170").Value = 1.975
$ws.Range("V170").Value = 1.825
